# Quarterly indexing esoteric bug-fix operation
#
# Column A holds quarter-index dates stored as the 1st day of the
# corresponding quarter-start month. The indexing was off: every date
# needs to be re-pointed to the 15th of the *following* month (i.e. the
# dates were effectively anchored 45 days too early). Walk every
# populated row in column A and correct the stored serial date in place,
# leaving every other cell/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)

    if ($cell.Value2 -eq $null) {
        continue
    }

    $origDate = [DateTime]::FromOADate($cell.Value2)
    $nextMonth = $origDate.AddMonths(1)
    $fixedDate = $nextMonth.AddDays(15 - $nextMonth.Day)

    $cell.Value = $fixedDate.ToOADate()
}
